$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = 0.03998220957160841
$ws.Range("J2").Value = 0.2646407987648076
$ws.Range("K2").Value = -0.1723740260921927
$ws.Range("L2").Value = 2.537203454282187

$ws.Range("I20").Value = -0.1582067356211394
$ws.Range("J20").Value = 0.3588061327005375
$ws.Range("K20").Value = 0.03761786831868028
$ws.Range("L20").Value = 2.061438656335365
